$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B45").Value = "Sundaram Fastners Ltd(TVS)-2"
$ws.Range("B45").Characters(1, 26).Font.Color = 6438691
$ws.Range("B45").Characters(27, 2).Font.Underline = $false
